$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $val)
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.405.53'
Set-TextValue 'E2' '  -1.90%  '
Set-TextValue 'D3' '1.797.50'
Set-TextValue 'E3' '  -1.71%  '
Set-TextValue 'D4' '1.006'
Set-TextValue 'E4' '  -0.11%  '
Set-TextValue 'D5' '1.006'
Set-TextValue 'E5' '  -0.06%  '
Set-TextValue 'D6' '307.85'
Set-TextValue 'E6' '  -1.10%  '
Set-TextValue 'D7' '0.4517'
Set-TextValue 'E7' '  -1.27%  '
Set-TextValue 'D8' '0.3600'
Set-TextValue 'E8' '  -2.64%  '
Set-TextValue 'E9' '  +0.00%  '
Set-TextValue 'D10' '0.07075'
Set-TextValue 'E10' '  -1.26%  '
Set-TextValue 'D11' '0.8859'
Set-TextValue 'E11' '  +1.09%  '
Set-TextValue 'D12' '0.07743'
Set-TextValue 'E12' '  -0.41%  '
Set-TextValue 'D13' '19.47'
Set-TextValue 'E13' '  -0.89%  '
Set-TextValue 'D14' '1.765.29'
Set-TextValue 'E14' '  -4.01%  '
Set-TextValue 'D15' '5.285'
Set-TextValue 'E15' '  -0.79%  '
Set-TextValue 'D16' '6.330'
Set-TextValue 'E16' '  -1.13%  '
Set-TextValue 'D17' '85.03'
Set-TextValue 'E17' '  -2.38%  '
Set-TextValue 'D18' '1.008'
Set-TextValue 'E18' '  -0.16%  '
Set-TextValue 'D19' '0.000008529'
Set-TextValue 'E19' '  -2.30%  '
Set-TextValue 'D20' '1.006'
Set-TextValue 'E20' '  -0.05%  '
Set-TextValue 'D21' '14.28'
Set-TextValue 'E21' '  -1.53%  '
Set-TextValue 'D22' '26.420.60'
Set-TextValue 'E22' '  -1.99%  '
Set-TextValue 'D23' '4.978'
Set-TextValue 'E23' '  -0.69%  '
Set-TextValue 'D24' '2.028.08'
Set-TextValue 'E24' '  -1.83%  '
Set-TextValue 'D25' '10.55'
Set-TextValue 'E25' '  +1.24%  '
Set-TextValue 'D26' '1.976'
Set-TextValue 'E26' '  -1.26%  '
Set-TextValue 'D27' '151.42'
Set-TextValue 'E27' '  +0.05%  '
Set-TextValue 'D28' '17.80'
Set-TextValue 'E28' '  -2.13%  '
Set-TextValue 'D29' '2.019'
Set-TextValue 'E29' '  +2.65%  '
Set-TextValue 'D30' '111.77'
Set-TextValue 'E30' '  -1.86%  '
Set-TextValue 'D31' '4.887'
Set-TextValue 'E31' '  -0.83%  '
Set-TextValue 'D32' '0.08686'
Set-TextValue 'E32' '  -1.38%  '
Set-TextValue 'D33' '3.082'
Set-TextValue 'E33' '  +1.68%  '
Set-TextValue 'B34' 'RenderToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D34' '2.740'
Set-TextValue 'E34' '  +6.78%  '
Set-TextValue 'B35' 'Filecoin'
Set-TextValue 'C35' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D35' '4.444'
Set-TextValue 'E35' '  -1.10%  '
Set-TextValue 'B36' 'ImmutableX'
Set-TextValue 'C36' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D36' '0.7227'
Set-TextValue 'E36' '  -3.74%  '
Set-TextValue 'D37' '1.107'
Set-TextValue 'E37' '  -2.56%  '
Set-TextValue 'D38' '1.005'
Set-TextValue 'E38' '  +0.09%  '
Set-TextValue 'D39' '1.067'
Set-TextValue 'E39' '  -2.05%  '
Set-TextValue 'D40' '0.01933'
Set-TextValue 'E40' '  -0.69%  '
Set-TextValue 'D41' '0.05101'
Set-TextValue 'E41' '  -0.91%  '
Set-TextValue 'D42' '2.867'
Set-TextValue 'E42' '  -1.44%  '
Set-TextValue 'B43' 'TheSandbox'
Set-TextValue 'C43' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D43' '0.5058'
Set-TextValue 'E43' '  +1.42%  '
Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '6.850'
Set-TextValue 'E44' '  -1.57%  '
Set-TextValue 'D45' '0.1517'
Set-TextValue 'E45' '  -5.20%  '
Set-TextValue 'D46' '8.013'
Set-TextValue 'E46' '  -3.85%  '
Set-TextValue 'D47' '1.007'
Set-TextValue 'E47' '  -0.03%  '
Set-TextValue 'D48' '0.4624'
Set-TextValue 'E48' '  -1.55%  '
Set-TextValue 'D49' '101.11'
Set-TextValue 'E49' '  -0.99%  '
Set-TextValue 'D50' '9.836'
Set-TextValue 'E50' '  -2.79%  '
Set-TextValue 'D51' '1.574'
Set-TextValue 'E51' '  -2.61%  '
